$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "66.311.72"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.54%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.323.16"
$c.ClearFormats()
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "587.42"
$c.ClearFormats()
$ws.Range("E5").Value = "  +2.54%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "184.18"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.37%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.649"
$c.ClearFormats()
$ws.Range("E7").Value = "  +7.86%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.41%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.83"
$c.ClearFormats()
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("E11").Value = "  +0.03%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "3.901.57"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  -3.18%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "66.331.23"
$c.ClearFormats()
$ws.Range("E14").Value = "  -0.65%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "26.21"
$c.ClearFormats()
$ws.Range("E15").Value = "  -3.15%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.326.02"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("E17").Value = "  -2.14%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "425.02"
$c.ClearFormats()
$ws.Range("E18").Value = "  -2.78%  "
$ws.Range("E19").Value = "  -2.40%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.21"
$c.ClearFormats()
$ws.Range("E20").Value = "  -2.69%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.40"
$c.ClearFormats()
$ws.Range("E21").Value = "  -2.52%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "71.86"
$c.ClearFormats()
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("E23").Value = "  +0.22%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.67"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.01%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.463.67"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("E26").Value = "  -0.56%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.203"
$c.ClearFormats()
$ws.Range("E27").Value = "  +6.62%  "
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -2.28%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "22.41"
$c.ClearFormats()
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("E35").Value = "  -3.01%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.18"
$c.ClearFormats()
$ws.Range("E36").Value = "  -3.94%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "159.82"
$c.ClearFormats()
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("E38").Value = "  -3.13%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.891.09"
$c.ClearFormats()
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("E40").Value = "  -1.74%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.44"
$c.ClearFormats()
$ws.Range("E41").Value = "  -5.03%  "
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -0.15%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0665"
$c.ClearFormats()
$ws.Range("E45").Value = "  -0.30%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "5.95"
$c.ClearFormats()
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("E47").Value = "  -1.82%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "23.27"
$c.ClearFormats()
$ws.Range("E48").Value = "  -5.24%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "314.75"
$c.ClearFormats()
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -0.50%  "
